# ==========================================================================
# 24 网络技术3班考勤.xlsx — apply "Add files via upload / 10.28" edit
#
#   1. Rename existing sheets: Sheet3 -> 考勤, Sheet1 -> 打扫卫生
#   2. Insert a brand-new sheet 实验报告 between them (lab-report roster for
#      the new 10.28 session)
#   3. Add a new "10.28" attendance column (R) to 考勤, with per-student
#      present/absent marks; the existing COUNTIF-based score formulas in
#      column C recalc automatically, except the one literal override cell
#      (C33) which is corrected by hand to match
#   4. Add a new "10.28" cleaning-duty column (E) to 打扫卫生
#   5. Make 实验报告 the active tab, matching the saved workbook view
# ==========================================================================

$wb = $excel.ActiveWorkbook

# --- 1. rename the two original sheets -----------------------------------
$wsAttendance = $wb.Worksheets.Item("Sheet3")
$wsAttendance.Name = "考勤"

$wsCleaning = $wb.Worksheets.Item("Sheet1")
$wsCleaning.Name = "打扫卫生"

# --- 2. insert the new lab-report sheet right after 考勤 -------------------
$wsLab = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsAttendance)
$wsLab.Name = "实验报告"

# Re-resolve the pre-existing sheet handles by name: inserting a new sheet
# can leave older worksheet object references stale (writes silently
# no-op), so fetch fresh ones before doing any further cell edits.
$wsAttendance = $wb.Worksheets.Item("考勤")
$wsCleaning = $wb.Worksheets.Item("打扫卫生")

# --- 3. 考勤: new "10.28" column (R) --------------------------------------
$wsAttendance.Range("R2").Value = 10.28

$attendanceMarks = [ordered]@{
    3  = "√"; 4  = "√"; 5  = "×"; 6  = "×"; 7  = "√"; 8  = "√"; 9  = "√"
    10 = "√"; 11 = "×"; 12 = "√"; 13 = "√"; 14 = "√"; 15 = "√"; 16 = "×"
    17 = "√"; 18 = "√"; 19 = "√"; 20 = "√"; 21 = "√"; 22 = "√"; 23 = "√"
    24 = "√"; 25 = "×"; 26 = "×"; 27 = "×"; 28 = "√"; 29 = "√"; 30 = "√"
    31 = "×"; 32 = "√"; 33 = "×"; 34 = "√"; 35 = "×"; 36 = "√"; 37 = "√"
    38 = "×"; 39 = "√"; 40 = "×"; 41 = "√"; 42 = "√"; 43 = "√"; 44 = "√"
    45 = "×"; 46 = "×"; 47 = "×"; 48 = "√"; 49 = "√"
}
foreach ($row in $attendanceMarks.Keys) {
    $wsAttendance.Cells.Item($row, 18).Value = $attendanceMarks[$row]
}

# C33 is a hand-typed literal (not the shared COUNTIF formula), so it needs
# an explicit correction to stay consistent with the new × in R33
$wsAttendance.Range("C33").Value = 95

# --- 4. 打扫卫生: new "10.28" column (E) -----------------------------------
$wsCleaning.Range("E2").Value = 10.28
$wsCleaning.Range("E3").Value = "江春阳"
$wsCleaning.Range("E4").Value = "孙蓓"
$wsCleaning.Range("E5").Value = "高明玉"
$wsCleaning.Range("E6").Value = "翟文铄"
$wsCleaning.Range("E8").Value = "刘佳鹏"
$wsCleaning.Range("E9").Value = "张婧"
$wsCleaning.Range("E10").Value = "王化坤"

# --- 5. 实验报告: header + roster for the 10.28 session --------------------
$wsLab.Range("A1").Value = "姓名"
$wsLab.Range("A1:A2").Merge()
$wsLab.Range("B2").Value = [DateTime]"2025-10-28"
$wsLab.Range("B2").NumberFormat = "m/d/yyyy"
$wsLab.Columns.Item(2).ColumnWidth = 15.625

$labRoster = @(
    @{row=3;  name="李奥";   status="×"}
    @{row=4;  name="官长皓"; status="×"}
    @{row=5;  name="刘柯纬"; status="×"}
    @{row=6;  name="张智献"; status="×"}
    @{row=7;  name="肖笛";   status="√"}
    @{row=8;  name="宁佳怡"; status="√"}
    @{row=9;  name="宋端祥"; status="√"}
    @{row=10; name="孙赫";   status="×"}
    @{row=11; name="高顼研"; status="×"}
    @{row=12; name="周欣慧"; status="√"}
    @{row=13; name="孔繁浩"; status="√"}
    @{row=14; name="王淑雨"; status="√"}
    @{row=15; name="刘建平"; status="√"}
    @{row=16; name="李志远"; status="×"}
    @{row=17; name="李雨欣"; status="√"}
    @{row=18; name="江春阳"; status="√"}
    @{row=19; name="黄一坤"; status="×"}
    @{row=20; name="朱锦涛"; status="√"}
    @{row=21; name="姜立敏"; status="√"}
    @{row=22; name="张在满"; status="√"}
    @{row=23; name="马圣涵"; status="√"}
    @{row=24; name="孙蓓";   status="×"}
    @{row=25; name="陈金腾"; status="×"}
    @{row=26; name="王安达"; status="×"}
    @{row=27; name="夹梦娅"; status="×"}
    @{row=28; name="石文凯"; status="√"}
    @{row=29; name="祁志一"; status="√"}
    @{row=30; name="梁亚伟"; status="×"}
    @{row=31; name="田光宁"; status="×"}
    @{row=32; name="高明玉"; status="×"}
    @{row=33; name="崔斐艳"; status="×"}
    @{row=34; name="翟文铄"; status="×"}
    @{row=35; name="刘佳鹏"; status="×"}
    @{row=36; name="张婧";   status="√"}
    @{row=37; name="王化坤"; status="√"}
    @{row=38; name="刘存铎"; status="×"}
    @{row=39; name="李有梁"; status="×"}
    @{row=40; name="陈时";   status="×"}
    @{row=41; name="于米朵"; status="√"}
    @{row=42; name="张初晨"; status="×"}
    @{row=43; name="刘倩惠"; status="√"}
    @{row=44; name="张顺";   status="×"}
    @{row=45; name="刘明义"; status="×"}
    @{row=46; name="郝江涛"; status="×"}
    @{row=47; name="侯明幸"; status="×"}
    @{row=48; name="葛欣宇"; status="√"}
    @{row=49; name="刘雨蒙"; status="√"}
)
foreach ($entry in $labRoster) {
    $wsLab.Cells.Item($entry.row, 1).Value = $entry.name
    $wsLab.Cells.Item($entry.row, 2).Value = $entry.status
}

# --- 6. make 实验报告 the active tab, matching the saved view -------------
$wsLab.Activate()
$wsLab.Range("B20").Select()

Write-Host "Done. Sheets:" ($wb.Worksheets | ForEach-Object { $_.Name })
